# Add a new "Chain Of Responsibility" pattern section right after the
# Prototype pattern's description paragraph.

$d = $word.ActiveDocument

# Locate the end of the Prototype pattern's description paragraph
# ("使用原型实例指定要创建对象的类型，通过复制这个原型来创建新对象。") -
# the new section is inserted right after it.
$rng = $d.Content
$found = $rng.Find.Execute(
    "使用原型实例指定要创建对象的类型，通过复制这个原型来创建新对象。",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the Prototype pattern description paragraph."
}

# Collapse to the end of the match (right before the following paragraph).
$rng.Collapse(0)

# Re-derive a plain Range at the same offsets. (A Range produced by
# Find+Collapse keeps extra internal state that confuses a multi-paragraph
# InsertXML -- it can clobber the preceding paragraph instead of purely
# inserting after it. A fresh Range built from the same Start/End avoids
# that.)
$insertionPoint = $d.Range($rng.Start, $rng.End)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$newParagraphsXml = @"
<w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">1. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>责任链（</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Chain Of Responsibility</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>）</w:t></w:r></w:p><w:p xmlns:w="$wNs"/><w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>意图</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>使多个对象都有机会处理请求，从而避免请求的发送者和接收者之间的耦合关系。将这些对象连成一条链，并沿着这条链发送该请求，直到有一个对象处理它为止。</w:t></w:r></w:p>
"@

$insertionPoint.InsertXML($newParagraphsXml)
